# Fix bug infinity loop: append the missing Room 505 row to the Rooms sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooms")

# Find the first empty row after the existing data (row 9, right after row 8).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "505"
$ws.Cells.Item($newRow, 2).Value = 19
$ws.Cells.Item($newRow, 3).Value = "Single"
$ws.Cells.Item($newRow, 4).Value = "No"
